# Automatische test-sync: 2025-08-08 20:00:50
#
# Adds a new log entry (row 3) to the "Logs" sheet and the corresponding
# aggregated "Overig" category (row 3) to the "Dashboard" sheet, extends the
# conditional formatting ranges to cover the new row, and widens the chart's
# series references so the new Dashboard row is plotted too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet - append the new e-mail log entry in row 3
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(3, 1).Value = "Zou jij dit even op kunnen pakken?"
$logs.Cells.Item(3, 2).Value = '"Testbedrijf 123 B.V." <admin@testbedrijf123.nl>'
$logs.Cells.Item(3, 3).Value = "Testmail #1: Zou jij dit even op kunnen pakken?`nTestbedrijf 123 B.V."
$logs.Cells.Item(3, 4).Value = "Overig"
$logs.Cells.Item(3, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@testbedrijf123.nl."
$logs.Cells.Item(3, 6).Value = "2025-08-08 19:59:52"
$logs.Cells.Item(3, 7).Value = "Ja"
$logs.Cells.Item(3, 8).Value = "Ja"
$logs.Cells.Item(3, 9).Value = "Nee"
$logs.Cells.Item(3, 10).Value = "Nee"

# Extend the conditional formatting ranges (D, G, H, I, J) so they cover the
# newly added row as well.
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $fcs = $logs.Range("$col`2").FormatConditions
    $newRange = $logs.Range("$col`2:$col`3")
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------------
# 2. Dashboard sheet - append the aggregated "Overig" count in row 3
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(3, 1).Value = "Overig"
$dashboard.Cells.Item(3, 2).Value = 1

# ---------------------------------------------------------------------------
# 3. Chart - widen the category/value series references to include row 3
# ---------------------------------------------------------------------------
$chartObj = $dashboard.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
